# Revert "Merge branch 'wrong-xlsform-col'"
#
# The only user-visible/content change in the target workbook is that the
# cursor / active-cell selection on the "survey" sheet moves from C2 back
# to B5 (the shared-strings table reshuffle seen in the raw OOXML diff is
# purely an artifact of how the original editor serialized the workbook
# and carries no cell-content difference - every cell keeps the exact same
# text it had before).

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Activate()
$survey.Range("B5").Select()
